$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy old row 9 data down to new row 10 (shift existing record down)
$ws.Range("A10").Value2 = $ws.Range("A9").Value2
$ws.Range("B10").Value2 = $ws.Range("B9").Value2
$ws.Range("C10").Value2 = $ws.Range("C9").Value2
$ws.Range("D10").Value2 = $ws.Range("D9").Value2
$ws.Range("D10").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("E10").Value2 = $ws.Range("E9").Value2
$ws.Range("F10").Value2 = $ws.Range("F9").Value2
$ws.Range("G10").Value2 = $ws.Range("G9").Value2
$ws.Range("H10").Value2 = $ws.Range("H9").Value2
$ws.Range("I10").Value2 = $ws.Range("I9").Value2
$ws.Range("J10").Value2 = $ws.Range("J9").Value2
$ws.Range("K10").Value2 = $ws.Range("K9").Value2
$ws.Range("L10").Value2 = $ws.Range("L9").Value2
$ws.Range("M10").Value2 = $ws.Range("M9").Value2
$ws.Range("N10").Value2 = $ws.Range("N9").Value2
$ws.Range("O10").Value2 = $ws.Range("O9").Value2
$ws.Range("P10").Value2 = $ws.Range("P9").Value2
$ws.Range("Q10").Value2 = $ws.Range("Q9").Value2
$ws.Range("R10").Value2 = $ws.Range("R9").Value2
$ws.Range("S10").Value2 = $ws.Range("S9").Value2
$ws.Range("T10").Value2 = $ws.Range("T9").Value2

# Now update row 9 with the new record values
$ws.Range("D9").Value2 = 44596
$ws.Range("M9").Value2 = 120
$ws.Range("N9").Value2 = 2500
$ws.Range("O9").Value2 = 2700
$ws.Range("P9").Value2 = 2600
$ws.Range("S9").Value2 = 1300
$ws.Range("T9").Value2 = 2
